$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bump post_datetime for rows 75-77 forward by exactly 1 hour (1/24 of a day) ---
foreach ($r in 75..77) {
    $cell = $ws.Cells.Item($r, 10)   # column J = post_datetime
    $orig = $cell.Value2
    $cell.Value2 = ($orig * 24 + 1) / 24
}

# --- Update the sheet's scroll position / active selection to match the new view ---
$win = $excel.ActiveWindow
$win.ScrollRow = 77
$win.ScrollColumn = 4
$null = $ws.Range("J76:J77").Select()

# --- Force a full recalculation so the volatile RANDBETWEEN-based K/L/M columns refresh ---
$null = $excel.CalculateFull()
